$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.91849999999999
$ws.Range("A9").Value = -20.16709999999998
$ws.Range("B11").Value = 5.548199999999997
$ws.Range("A18").Value = -23.1552
$ws.Range("A20").Value = -22.10570000000002
